$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Make room for the new "StrikeCraft" weapons block by inserting 4 blank rows
# right after the existing row 55 (the old "Torpedo Type" section header).
# This pushes the old rows 56-58 down to 60-62 and Excel auto-repoints the
# formulas that referenced those rows (e.g. D56 -> D60).
# ---------------------------------------------------------------------------
$ws.Rows("56:59").Insert()

# ---------------------------------------------------------------------------
# Row 55 becomes the header for the new StrikeCraft block (it mirrors row 1):
# "Size","Weapon","Ammo","SD","AP","AD","CD","HD","HT","Systems hit",
# "Blast Radius","Weapon effect size". Columns F-L/A keep their old values.
# ---------------------------------------------------------------------------
$ws.Range("B55").Value = "Size"
$ws.Range("C55").Value = "Weapon"
$ws.Range("E55").Value = "Ammo"
$ws.Range("M55").Value = "Blast Radius"
$ws.Range("N55").Value = "Weapon effect size"

# ---------------------------------------------------------------------------
# Row 56: StrikeCraft / FighterCannon
# ---------------------------------------------------------------------------
$ws.Range("A56").Value = 3
$ws.Range("B56").Value = "StrikeCraft"
$ws.Range("C56").Value = "FighterCannon"
$ws.Range("D56").Value = 0.75
$ws.Range("E56").Value = "KineticPenetrator"
$ws.Range("F56").Formula = "=1*D56*P20"
$ws.Range("G56").Value = 175
$ws.Range("H56").Formula = "=1*D56*Q20"
$ws.Range("I56").Formula = "=1*D56*R20"
$ws.Range("J56").Formula = "=1*D56*S20"
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 1
$ws.Range("M56").Value = 0
$ws.Range("N56").Value = 0.3

# ---------------------------------------------------------------------------
# Row 57: StrikeCraft / FighterAutoannon
# ---------------------------------------------------------------------------
$ws.Range("A57").Value = 3
$ws.Range("B57").Value = "StrikeCraft"
$ws.Range("C57").Value = "FighterAutoannon"
$ws.Range("D57").Value = 0.5
$ws.Range("E57").Value = "KineticPenetrator"
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = 0
$ws.Range("N57").Value = 0.2

# Row 58 is left blank on purpose (matches the blank separator row pattern
# used throughout the sheet between sections).

# ---------------------------------------------------------------------------
# Row 59: re-create the "Torpedo Type" section header that used to live in
# row 55, now that row 55 has been repurposed.
# ---------------------------------------------------------------------------
$ws.Range("A59").Value = "#"
$ws.Range("B59").Value = "Torpedo Type"
$ws.Range("C59").Value = "Spread Size"
$ws.Range("E59").Value = "Range"
$ws.Range("F59").Value = "SD"
$ws.Range("G59").Value = "AP"
$ws.Range("H59").Value = "AD"
$ws.Range("I59").Value = "CD"
$ws.Range("J59").Value = "HD"
$ws.Range("K59").Value = "HT"
$ws.Range("L59").Value = "Systems hit"
$ws.Range("M59").Value = "Weapon effect size"
$ws.Range("N59").Value = "Projectile Size"

# ---------------------------------------------------------------------------
# Column C needs to be a bit wider to fit "FighterAutoannon".
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 16.592447916666668

# ---------------------------------------------------------------------------
# Update the view: scrolled down a bit further, new active cell.
# ---------------------------------------------------------------------------
$excel.Goto($ws.Range("A40"), $true)
$ws.Range("M56").Select()
